# Fill in row 10 of Sheet1 with the new match data (Alverca 0 - 4 Gil Vicente)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = "24/10/2025"
$ws.Range("B10").Value = "Alverca"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = "Gil Vicente"
$ws.Range("F10").Value = "W"
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1.21
$ws.Range("L10").Value = 1.56
$ws.Range("M10").Value = 12
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 5
$ws.Range("P10").Value = 2
